# Applies the tc_p120v.docx edit described by the commit:
# "rewrite of rclone and preprocessing script; update xml"
#
# 1) Fix typo "il en ya" -> "il en y a"
# 2) Drop a stray <add>...</add> TEI wrapper around a lone "e"
# 3) Drop a stray <add>De</add> TEI wrapper, merging its text into
#    the following run ("De" + " cela " -> "De cela ")
# 4) Drop a stray "." run before a "</add>" marker
# 5) Drop a manual line break after "<rub>"
# 6) Move the "</rub>" marker to before " le second foeillet suivant"
#    (instead of after it), matching the TEI markup correction

$d = $word.ActiveDocument

# --- 1) typo fix -----------------------------------------------------
$null = $d.Content.Find.Execute(
    "il en ya les rend", $true, $false, $false, $false, $false,
    $true, 1, $false, "il en y a les rend", 2)

# --- 2) drop <add>...</add> wrapper around the lone "e" --------------
# "d<del>u</del><add>e</add> lantimoyne" -> "d<del>u</del>e lantimoyne"
$t = $d.Content.Text
$idx = $t.IndexOf("<del>u</del><add>e</add> lantimoyne")
$openStart = $idx + ("<del>u</del>").Length
$openEnd = $openStart + ("<add>").Length
$closeStart = $openEnd + ("e").Length
$closeEnd = $closeStart + ("</add>").Length
# delete the closing tag first (rightmost) so offsets of the opening tag
# stay valid
$d.Range($closeStart, $closeEnd).Delete()
$d.Range($openStart, $openEnd).Delete()

# --- 3) drop <add>...</add> wrapper around "De" -----------------------
# "<del>Cela</del> <add>De</add> cela " -> "<del>Cela</del> De cela "
$t = $d.Content.Text
$idx = $t.IndexOf("<del>Cela</del> <add>De</add> cela ")
$openStart = $idx + ("<del>Cela</del> ").Length
$openEnd = $openStart + ("<add>").Length
$closeStart = $openEnd + ("De").Length
$closeEnd = $closeStart + ("</add>").Length
$d.Range($closeStart, $closeEnd).Delete()
$d.Range($openStart, $openEnd).Delete()

# --- 4) drop the stray "." run right before the closing </add> -------
# "antimoyne</m> &amp;.</add>" -> "antimoyne</m> &amp;</add>"
$t = $d.Content.Text
$idx = $t.IndexOf("antimoyne</m> &amp;.</add>")
$dotStart = $idx + ("antimoyne</m> &amp;").Length
$dotEnd = $dotStart + (".").Length
$d.Range($dotStart, $dotEnd).Delete()

# --- 5) drop the manual line break right after "<rub>" ----------------
$t = $d.Content.Text
$idx = $t.IndexOf("<rub>")
$breakStart = $idx + ("<rub>").Length
$breakEnd = $breakStart + 1
$d.Range($breakStart, $breakEnd).Delete()

# --- 6) move "</rub>" from after "suivant" to before " le second ------
#        foeillet suivant" -------------------------------------------
# Capture the exact run formatting of the existing "</rub>" text (a
# Courier New / blue TEI-tag run identical to the "<rub>" one) so the
# freshly inserted copy keeps the same rPr, then delete the old one.
$t = $d.Content.Text
$idx = $t.IndexOf("<rub>")
$tagStart = $idx
$tagEnd = $tagStart + 1   # just the leading "<" run, same formatting as </rub>
$srcFormatted = $d.Range($tagStart, $tagEnd).FormattedText

$closeIdx = $t.IndexOf("</rub>")
$closeEnd = $closeIdx + ("</rub>").Length
$d.Range($closeIdx, $closeEnd).Delete()

$t = $d.Content.Text
$insertIdx = $t.IndexOf(" le second foeillet suivant")
$destRange = $d.Range($insertIdx, $insertIdx)
$destRange.FormattedText = $srcFormatted
$newRange = $d.Range($insertIdx, $insertIdx + 1)
$newRange.Text = "</rub>"

Write-Output "done"
